$wb = $excel.ActiveWorkbook

# --- Sheet1 (peak_table): update w_height (N) / m_height (O) ---
$ws1 = $wb.Worksheets.Item("peak_table")
$ws1.Cells.Item(2, 14).Value = 800
$ws1.Cells.Item(2, 15).Value = 300
$ws1.Cells.Item(3, 14).Value = 1000
$ws1.Cells.Item(3, 15).Value = 900
$ws1.Cells.Item(4, 14).Value = 400
$ws1.Cells.Item(4, 15).Value = 1000
$ws1.Cells.Item(6, 14).Value = 700
$ws1.Cells.Item(6, 15).Value = 1000
$ws1.Cells.Item(12, 14).Value = 600
$ws1.Cells.Item(12, 15).Value = 1000
$ws1.Cells.Item(14, 14).Value = 700
$ws1.Cells.Item(14, 15).Value = 1000

# --- Sheet2 (allele_table): update min_height (K), is_detected (M), peak (N), size (O), height (P), status (Q), message (R) ---
$ws2 = $wb.Worksheets.Item("allele_table")
$ws2.Cells.Item(2, 11).Value = 800
$ws2.Cells.Item(2, 13).Value = $True
$ws2.Cells.Item(2, 14).Value = 37
$ws2.Cells.Item(2, 15).Value = 29.15
$ws2.Cells.Item(2, 16).Value = 886
$ws2.Cells.Item(2, 17).Value = "ok"
$ws2.Cells.Item(2, 18).Value = ""
$ws2.Cells.Item(3, 11).Value = 300
$ws2.Cells.Item(3, 13).Value = $True
$ws2.Cells.Item(3, 14).Value = 58
$ws2.Cells.Item(3, 15).Value = 33.43
$ws2.Cells.Item(3, 16).Value = 428
$ws2.Cells.Item(3, 17).Value = "ok"
$ws2.Cells.Item(3, 18).Value = ""
$ws2.Cells.Item(5, 11).Value = 900
$ws2.Cells.Item(5, 13).Value = $True
$ws2.Cells.Item(5, 14).Value = 38
$ws2.Cells.Item(5, 15).Value = 35.15
$ws2.Cells.Item(5, 16).Value = 943
$ws2.Cells.Item(5, 17).Value = "ok"
$ws2.Cells.Item(5, 18).Value = ""
$ws2.Cells.Item(6, 11).Value = 400
$ws2.Cells.Item(6, 13).Value = $True
$ws2.Cells.Item(6, 14).Value = 17
$ws2.Cells.Item(6, 15).Value = 38.89
$ws2.Cells.Item(6, 16).Value = 484
$ws2.Cells.Item(6, 17).Value = "ok"
$ws2.Cells.Item(6, 18).Value = ""
$ws2.Cells.Item(10, 11).Value = 700
$ws2.Cells.Item(10, 13).Value = $True
$ws2.Cells.Item(10, 14).Value = 18
$ws2.Cells.Item(10, 15).Value = 46.9
$ws2.Cells.Item(10, 16).Value = 914
$ws2.Cells.Item(10, 17).Value = "ok"
$ws2.Cells.Item(10, 18).Value = ""
$ws2.Cells.Item(22, 11).Value = 600
$ws2.Cells.Item(22, 13).Value = $True
$ws2.Cells.Item(22, 14).Value = 42
$ws2.Cells.Item(22, 15).Value = 30.66
$ws2.Cells.Item(22, 16).Value = 655
$ws2.Cells.Item(22, 17).Value = "ok"
$ws2.Cells.Item(22, 18).Value = ""
$ws2.Cells.Item(26, 11).Value = 700
$ws2.Cells.Item(26, 13).Value = $True
$ws2.Cells.Item(26, 14).Value = 26
$ws2.Cells.Item(26, 15).Value = 38.83
$ws2.Cells.Item(26, 16).Value = 884
$ws2.Cells.Item(26, 17).Value = "ok"
$ws2.Cells.Item(26, 18).Value = ""

# --- Sheet3 (marker_table): update genotype (G) / phenotype (H) ---
$ws3 = $wb.Worksheets.Item("marker_table")
$ws3.Cells.Item(2, 7).Value = "GA"
$ws3.Cells.Item(2, 8).Value = "heterozygous"
$ws3.Cells.Item(3, 7).Value = "T"
$ws3.Cells.Item(3, 8).Value = "homozygous mutant"
$ws3.Cells.Item(4, 7).Value = "T"
$ws3.Cells.Item(4, 8).Value = "wildtype"
$ws3.Cells.Item(5, 7).Value = "G"
$ws3.Cells.Item(5, 8).Value = "wildtype"
$ws3.Cells.Item(6, 7).Value = "G"
$ws3.Cells.Item(6, 8).Value = "wildtype"
$ws3.Cells.Item(7, 7).Value = "G"
$ws3.Cells.Item(7, 8).Value = "wildtype"
$ws3.Cells.Item(8, 7).Value = "G"
$ws3.Cells.Item(8, 8).Value = "wildtype"
$ws3.Cells.Item(9, 7).Value = "C"
$ws3.Cells.Item(9, 8).Value = "wildtype"
$ws3.Cells.Item(10, 7).Value = "G"
$ws3.Cells.Item(10, 8).Value = "wildtype"
$ws3.Cells.Item(12, 7).Value = "G"
$ws3.Cells.Item(12, 8).Value = "wildtype"
$ws3.Cells.Item(13, 7).Value = "A"
$ws3.Cells.Item(13, 8).Value = "wildtype"
$ws3.Cells.Item(14, 7).Value = "C"
$ws3.Cells.Item(14, 8).Value = "wildtype"
$ws3.Cells.Item(15, 7).Value = "A"
$ws3.Cells.Item(15, 8).Value = "wildtype"
$ws3.Cells.Item(16, 7).Value = "T"
$ws3.Cells.Item(16, 8).Value = "wildtype"
$ws3.Cells.Item(17, 7).Value = "G"
$ws3.Cells.Item(17, 8).Value = "wildtype"
$ws3.Cells.Item(18, 7).Value = "G"
$ws3.Cells.Item(18, 8).Value = "wildtype"

# --- Leave the view focused back on peak_table, matching the reviewed cell ---
$ws1.Activate()
$ws1.Range("N6").Select()
